# Update the Gdnf-Gfra1 LR-pairs sheet with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster=MuSCs, Target cluster=ECs) - receptor-side and edge-weight
# values recomputed with new TPM data.
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01036033333333333
$ws.Range("N2").Value = 0.031081
$ws.Range("O2").Value = 0.0003369947480386084
$ws.Range("P2").Value = 0.0003369947480386084
$ws.Range("Q2").Value = 0.007230652759000001
$ws.Range("R2").Value = 0.06507587483100001
$ws.Range("S2").Value = 0.0003369947480386084
$ws.Range("T2").Value = 0.0003369947480386084

# Row 3 (Sending cluster=MuSCs, Target cluster=FAPs) - derived specificity values
# recomputed with new TPM data.
$ws.Range("O3").Value = 0.8439700329797517
$ws.Range("P3").Value = 0.8439700329797518
$ws.Range("S3").Value = 0.8439700329797517
$ws.Range("T3").Value = 0.8439700329797518

# Row 4 (Sending cluster=MuSCs, Target cluster=MuSCs) - derived specificity values
# recomputed with new TPM data.
$ws.Range("O4").Value = 0.1556929722722096
$ws.Range("P4").Value = 0.1556929722722096
$ws.Range("S4").Value = 0.1556929722722096
$ws.Range("T4").Value = 0.1556929722722096
